$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("shailesh@prac.to",  "['Dermatologist']", 5,  9,  "9087654321", "N"),
    @("uthappa@prac.to",   "['Cardiologist']",  10, 8,  "9087654321", "N"),
    @("sheela@prac.to",    "['Dermatologist']", 4,  8,  "9087654321", "N"),
    @("uthmini@prac.to",   "['Cardiologist']",  9,  7,  "9087654321", "N"),
    @("raj@prac.to",       "['General']",       20, 10, "9087654321", "N"),
    @("rajani@prac.to",    "['General']",       19, 9,  "9087654321", "N"),
    @("pulasthya@prac.to", "['Surgeon']",       1,  1,  "9087654321", "N"),
    @("Pulsathi@prac.to",  "['Surgeon']",       1,  1,  "9087654321", "N")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = "'" + $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row = $row + 1
}
